# "remove column from alcohol data"
# The alcohol measurement sheet has a duplicate/erroneous column (M) sitting
# immediately before the real data column (N). Delete the entire column M so
# every cell shifts one column to the left (N -> M, O -> N, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M:M").EntireColumn.Delete()

# Leave the selection on the (now left-shifted) former "N" column.
$ws.Range("M1").Select()
